$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 392, shifting existing rows 392:485 down to 393:486
$ws.Rows.Item(392).Insert()

# Populate the new row 392 with the new weekly data point
$ws.Range("A392").Value = 10
$ws.Range("B392").Value = "Vega Modelo de Temuco"
$ws.Range("C392").Value = "La Araucanía"
$ws.Range("D392").Value = 44736
$ws.Range("E392").Value = 9
$ws.Range("F392").Value = "Fruta"
$ws.Range("G392").Value = 100108
$ws.Range("H392").Value = "Tropicales y subtropicales"
$ws.Range("I392").Value = 100108005
$ws.Range("J392").Value = "Piña"
$ws.Range("K392").Value = "Caramelo"
$ws.Range("L392").Value = "Primera"
$ws.Range("M392").Value = 170
$ws.Range("N392").Value = 19000
$ws.Range("O392").Value = 20000
$ws.Range("P392").Value = 19529
$ws.Range("Q392").Value = "$/caja 12 unidades"
$ws.Range("R392").Value = "Ecuador"
$ws.Range("S392").Value = 1627
$ws.Range("T392").Value = 12

# Ensure the date cell keeps the same date/time number format as the rest of column D
$ws.Range("D392").NumberFormat = $ws.Range("D393").NumberFormat
